# "For test uploads scoot up points"
# Nudge the latitude column (P) north by +1 degree for every sample row,
# and correct the H8 flag (back to TRUE) on the GPS sample that was
# marked wrong while the points were being re-checked.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column P (latitude) holds values in rows 2-10; scoot each one up by 1.
for ($r = 2; $r -le 10; $r++) {
    $cell = $ws.Cells.Item($r, 16)   # column P
    $cell.Value2 = $cell.Value2 + 1
}

# Row 8's "present" flag (column H) flips back to TRUE.
$ws.Cells.Item(8, 8).Value2 = $true

# Leave the selection on the touched column.
$null = $ws.Range("P2:P10").Select()
